$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 48: duplicate of row 47 ("2021-01-09" / "10 Jan -- 16 Jan 2021" / 126.93 / KNN)
$ws.Range("A47:K47").Copy()
$ws.Range("A48:K48").PasteSpecial(-4163)

# Row 49: new forecast week ("2021-01-09" / "17 Jan -- 23 Jan 2021" / 125.24 / KNN)
$ws.Range("A47:F47").Copy()
$ws.Range("A49:F49").PasteSpecial(-4163)
$ws.Range("B49").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D49").Value = 125.24
